$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.912.24"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.230.96"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'302.93"
$ws.Range("E5").Value = "  -4.90%  "
$ws.Range("E6").Value = "  -7.29%  "
$ws.Range("D7").Value = "'0.566"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("E9").Value = "  -6.91%  "
$ws.Range("D10").Value = "'34.29"
$ws.Range("E10").Value = "  -7.86%  "
$ws.Range("D11").Value = "'0.0801"
$ws.Range("E11").Value = "  -3.62%  "
$ws.Range("E12").Value = "  -6.93%  "
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("D14").Value = "2.570.90"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").Value = "2.262.45"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "'0.808"
$ws.Range("E16").Value = "  -5.94%  "
$ws.Range("D17").Value = "'13.29"
$ws.Range("E17").Value = "  -8.44%  "
$ws.Range("D18").Value = "43.752.15"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "0.0₃0948"
$ws.Range("E19").Value = "  -3.87%  "
$ws.Range("D20").Value = "'11.99"
$ws.Range("E20").Value = "  -12.02%  "
$ws.Range("E21").Value = "  -6.42%  "
$ws.Range("D22").Value = "'64.32"
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").Value = "'235.54"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  -7.67%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -8.25%  "
$ws.Range("D27").Value = "'9.75"
$ws.Range("E27").Value = "  -4.43%  "
$ws.Range("E28").Value = "  -2.87%  "
$ws.Range("D29").Value = "'35.95"
$ws.Range("E29").Value = "  -3.46%  "
$ws.Range("D30").Value = "'19.91"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("D31").Value = "'5.83"
$ws.Range("E31").Value = "  -6.36%  "
$ws.Range("D32").Value = "'152.00"
$ws.Range("E32").Value = "  -4.42%  "
$ws.Range("D33").Value = "'0.0799"
$ws.Range("E33").Value = "  -6.32%  "
$ws.Range("D34").Value = "'2.63"
$ws.Range("E34").Value = "  -2.14%  "
$ws.Range("D35").Value = "'3.22"
$ws.Range("E35").Value = "  +4.56%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.108"
$ws.Range("E36").Value = "  -6.60%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "'0.117"
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("E38").Value = "  -10.30%  "
$ws.Range("D39").Value = "'14.57"
$ws.Range("E39").Value = "  -10.73%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").Value = "'3.29"
$ws.Range("E40").Value = "  -11.95%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'3.74"
$ws.Range("E41").Value = "  -11.13%  "
$ws.Range("D42").Value = "'0.0294"
$ws.Range("E42").Value = "  -6.67%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").Value = "1.716.36"
$ws.Range("E44").Value = "  -5.17%  "
$ws.Range("D45").Value = "'83.90"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'98.85"
$ws.Range("E46").Value = "  -5.48%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.183"
$ws.Range("E47").Value = "  -7.63%  "
$ws.Range("D48").Value = "'4.82"
$ws.Range("E48").Value = "  -7.82%  "
$ws.Range("D49").Value = "'14.30"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'7.97"
$ws.Range("E50").Value = "  -4.33%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'67.76"
$ws.Range("E51").Value = "  -10.85%  "
